$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A/B/C shift to B/C/D.
$ws.Columns("A").Insert()

# New column A header: copy the header formatting from the (shifted) former
# header cell B1, then overwrite its value with the new header text.
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "testCaseName"
$ws.Range("A2").Value = "TC_001"
$ws.Range("A3").Value = "TC_002"

# Re-point the hyperlinks: they used to live on column A (now column B).
# Stash B2's pristine formatting (the "Hyperlink" style) in a scratch cell so
# we can restore it after Hyperlinks.Add mutates the cell style.
$ws.Range("B2").Copy($ws.Range("Z1"))
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:automation@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:automation@gmail.com")
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# Column widths: new col A is narrower (test case id). Columns B/C already
# kept their original widths from the column-insert shift, so only A needs
# to be set. (13.14 lands exactly on the raw width=14 quantization step.)
$ws.Columns("A").ColumnWidth = 13.14

# Selection moves to C13.
$ws.Range("C13").Select()

# Add print/page setup (paper size 9 = A4, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
